$wb = $excel.ActiveWorkbook
# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2299.4
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H132").Value = 1302.2693
$ws.Range("I132").Value = 1189.5217
$ws.Range("J132").Value = 2166.6667
$ws.Range("K132").Value = 3568.5651
$ws.Range("L132").Value = 6500.000100000001
$ws.Range("M132").Value = -1038.5651
$ws.Range("N132").Value = -11560.0001
$ws.Range("H135").Value = 706.6
$ws.Range("I135").Value = 714.44446
$ws.Range("J135").Value = 636
$ws.Range("K135").Value = 6430.00014
$ws.Range("L135").Value = 5724
$ws.Range("M135").Value = -3895.00014
$ws.Range("N135").Value = -10794
$ws.Range("H137").Value = 1399
$ws.Range("I137").Value = 1343.1111
$ws.Range("J137").Value = 1524.75
$ws.Range("K137").Value = 4029.3333
$ws.Range("L137").Value = 4574.25
$ws.Range("M137").Value = -1479.3333
$ws.Range("N137").Value = -9674.25
$ws.Range("H138").Value = 3210.8696
$ws.Range("J138").Value = 2945.7778
$ws.Range("L138").Value = 8837.3334
$ws.Range("N138").Value = -19117.3334

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 152.5
$ws.Range("I5").Value = 152.5
$ws.Range("K5").Value = 152.5
$ws.Range("M5").Value = -40.5
$ws.Range("H45").Value = 1705.65
$ws.Range("I45").Value = 1583.909
$ws.Range("K45").Value = 1583.909
$ws.Range("M45").Value = -1206.909
$ws.Range("H61").Value = 4252.3477
$ws.Range("I61").Value = 2739.1667
$ws.Range("K61").Value = 2739.1667
$ws.Range("M61").Value = -2527.1667
$ws.Range("H74").Value = 1390.8286
$ws.Range("I74").Value = 873.5484
$ws.Range("K74").Value = 873.5484
$ws.Range("M74").Value = 0.4515999999999849
$ws.Range("H77").Value = 1390.8286
$ws.Range("I77").Value = 873.5484
$ws.Range("K77").Value = 4367.742
$ws.Range("M77").Value = 0.2579999999998108
$ws.Range("H97").Value = 513.1818
$ws.Range("I97").Value = 464.42856
$ws.Range("K97").Value = 464.42856
$ws.Range("M97").Value = 31.57144
$ws.Range("H122").Value = 1300
$ws.Range("I122").Value = 1300
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 3900
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -1450
$ws.Range("N122").Value = -8800
$ws.Range("H136").Value = 4252.3477
$ws.Range("I136").Value = 2739.1667
$ws.Range("K136").Value = 8217.500100000001
$ws.Range("M136").Value = -5667.500100000001
$ws.Range("H138").Value = 71214.5
$ws.Range("J138").Value = 71214.5
$ws.Range("L138").Value = 71214.5
$ws.Range("N138").Value = -81494.5

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 152.5
$ws.Range("I4").Value = 152.5
$ws.Range("K4").Value = 152.5
$ws.Range("M4").Value = -37.5
$ws.Range("H86").Value = 102294.05
$ws.Range("I86").Value = 2170.1428
$ws.Range("K86").Value = 2170.1428
$ws.Range("M86").Value = -1047.1428
$ws.Range("H89").Value = 102294.05
$ws.Range("I89").Value = 2170.1428
$ws.Range("K89").Value = 10850.714
$ws.Range("M89").Value = -5234.714
$ws.Range("H94").Value = 294.1875
$ws.Range("I94").Value = 298.8387
$ws.Range("K94").Value = 298.8387
$ws.Range("M94").Value = 152.1613
$ws.Range("H99").Value = 1826.125
$ws.Range("I99").Value = 1521.8
$ws.Range("K99").Value = 1521.8
$ws.Range("M99").Value = -23.79999999999995

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9566.667
$ws.Range("J4").Value = 9566.667
$ws.Range("L4").Value = 9566.667
$ws.Range("N4").Value = -9790.667
$ws.Range("H31").Value = 3015.3
$ws.Range("I31").Value = 2411.1035
$ws.Range("K31").Value = 2411.1035
$ws.Range("M31").Value = -2116.1035
$ws.Range("H34").Value = 3015.3
$ws.Range("I34").Value = 2411.1035
$ws.Range("K34").Value = 2411.1035
$ws.Range("M34").Value = -2209.1035
$ws.Range("H132").Value = 2018.1428
$ws.Range("I132").Value = 1067.3334
$ws.Range("J132").Value = 3444.3572
$ws.Range("K132").Value = 3202.0002
$ws.Range("L132").Value = 10333.0716
$ws.Range("M132").Value = -672.0002
$ws.Range("N132").Value = -15393.0716
$ws.Range("H134").Value = 980.53845
$ws.Range("I134").Value = 845.73914
$ws.Range("K134").Value = 2537.21742
$ws.Range("M134").Value = -2.217419999999947

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 166.57143
$ws.Range("I4").Value = 166.57143
$ws.Range("K4").Value = 499.71429
$ws.Range("M4").Value = -387.71429
$ws.Range("H131").Value = 764.89
$ws.Range("I131").Value = 440.625
$ws.Range("J131").Value = 793.087
$ws.Range("K131").Value = 1321.875
$ws.Range("L131").Value = 2379.261
$ws.Range("M131").Value = 3718.125
$ws.Range("N131").Value = -12459.261

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 118
$ws.Range("I2").Value = 192.66667
$ws.Range("J2").Value = 68.22222
$ws.Range("K2").Value = 192.66667
$ws.Range("L2").Value = 68.22222
$ws.Range("M2").Value = -79.66667000000001
$ws.Range("N2").Value = -294.22222
$ws.Range("H80").Value = 2947.1538
$ws.Range("I80").Value = 2842.375
$ws.Range("K80").Value = 2842.375
$ws.Range("M80").Value = -1844.375
$ws.Range("H83").Value = 2947.1538
$ws.Range("I83").Value = 2842.375
$ws.Range("K83").Value = 14211.875
$ws.Range("M83").Value = -9219.875
$ws.Range("H97").Value = 877.8261
$ws.Range("I97").Value = 827.93335
$ws.Range("K97").Value = 827.93335
$ws.Range("M97").Value = -331.93335
$ws.Range("H132").Value = 3183.4614
$ws.Range("I132").Value = 2541.5715
$ws.Range("K132").Value = 7624.7145
$ws.Range("M132").Value = -5094.7145

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 432727.28
$ws.Range("J2").Value = 130000
$ws.Range("L2").Value = 130000
$ws.Range("N2").Value = -130224
$ws.Range("H61").Value = 1757.4667
$ws.Range("I61").Value = 1689.8182
$ws.Range("J61").Value = 1943.5
$ws.Range("K61").Value = 1689.8182
$ws.Range("L61").Value = 1943.5
$ws.Range("M61").Value = -1487.8182
$ws.Range("N61").Value = -2347.5
$ws.Range("H113").Value = 1757.4667
$ws.Range("I113").Value = 1689.8182
$ws.Range("J113").Value = 1943.5
$ws.Range("K113").Value = 1689.8182
$ws.Range("L113").Value = 1943.5
$ws.Range("M113").Value = 480.1818000000001
$ws.Range("N113").Value = -6283.5

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 71448780
$ws.Range("I2").Value = 1000000000
$ws.Range("J2").Value = 21769.924
$ws.Range("K2").Value = 1000000000
$ws.Range("L2").Value = 21769.924
$ws.Range("M2").Value = -999999888
$ws.Range("N2").Value = -21993.924
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H135").Value = 71278.1
$ws.Range("J135").Value = 71278.1
$ws.Range("L135").Value = 71278.1
$ws.Range("N135").Value = -81418.1
$ws.Range("H136").Value = 2655.7666
$ws.Range("I136").Value = 2870.9333
$ws.Range("J136").Value = 2440.6
$ws.Range("K136").Value = 8612.7999
$ws.Range("L136").Value = 7321.799999999999
$ws.Range("M136").Value = -6062.7999
$ws.Range("N136").Value = -12421.8
